$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LISTE des sites web OFF")

# 1. Remove the two rows that are no longer "OFF" (marne.gouv.fr, martinique.gouv.fr)
$ws.Rows("363:364").Delete()

# 2. Update statuses for domains that are now redirecting instead of fully down
$ws.Range("B71").Value = "Redirection"
$ws.Range("B72").Value = "Redirection"

# permisdeconduire.gouv.fr was originally at row 407, now (after the 2-row delete) at row 405
$ws.Range("B405").Value = "Redirection"

# 24octobre.gouv.fr / 3939.gouv.fr move from "Possiblement OFF" to "Redirect"
$ws.Range("B593").Value = "Redirect"
$ws.Range("B594").Value = "Redirect"

# 3. View / window adjustments
$wb.Windows.Item(1).Width = 13960
$ws.Range("A485").Select()
